# Atualizado por script em 26-11-2023 20:30
# Append 3 new match rows (rows 30-32) to the Gibraltar National League sheet,
# matching the existing table's layout/formatting exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles) of the last existing data row (29) down onto
# the three new rows so column A keeps its bold/bordered/centered style (s=1)
# and column E keeps its date-time number format (s=2), same as every other
# data row in the sheet.
$ws.Range("A29:V29").Copy()
$ws.Range("A30:V32").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 30: Lions Gibraltar 0 - 0 Mons Calpe
$ws.Range("A30").Value = 29
$ws.Range("B30").Value = "gibraltar"
$ws.Range("C30").Value = "national-league"
$ws.Range("D30").Value = "2023-2024"
$ws.Range("E30").Value = 45254.875
$ws.Range("F30").Value = "Lions Gibraltar"
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = "Mons Calpe"
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 13.1
$ws.Range("K30").Value = "24/11/2023 10:32"
$ws.Range("L30").Value = 9.98
$ws.Range("M30").Value = "24/11/2023 20:20"
$ws.Range("N30").Value = 12.98
$ws.Range("O30").Value = "24/11/2023 10:32"
$ws.Range("P30").Value = 8.619999999999999
$ws.Range("Q30").Value = "24/11/2023 20:20"
$ws.Range("R30").Value = 1.06
$ws.Range("S30").Value = "24/11/2023 10:32"
$ws.Range("T30").Value = 1.13
$ws.Range("U30").Value = "24/11/2023 20:20"
$ws.Range("V30").Value = "https://www.betexplorer.com/football/gibraltar/national-league/lions-gibraltar-mons-calpe/Qcpa5lNf/"

# Row 31: College 1975 FC 0 - 5 Europa Point
$ws.Range("A31").Value = 30
$ws.Range("B31").Value = "gibraltar"
$ws.Range("C31").Value = "national-league"
$ws.Range("D31").Value = "2023-2024"
$ws.Range("E31").Value = 45255.6875
$ws.Range("F31").Value = "College 1975 FC"
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = "Europa Point"
$ws.Range("I31").Value = 5
$ws.Range("J31").Value = 2.23
$ws.Range("K31").Value = "25/11/2023 12:46"
$ws.Range("L31").Value = 2.32
$ws.Range("M31").Value = "25/11/2023 16:15"
$ws.Range("N31").Value = 3.53
$ws.Range("O31").Value = "25/11/2023 12:46"
$ws.Range("P31").Value = 3.72
$ws.Range("Q31").Value = "25/11/2023 16:15"
$ws.Range("R31").Value = 2.57
$ws.Range("S31").Value = "25/11/2023 12:46"
$ws.Range("T31").Value = 2.47
$ws.Range("U31").Value = "25/11/2023 16:15"
$ws.Range("V31").Value = "https://www.betexplorer.com/football/gibraltar/national-league/college-1975-europa-point/Or3248x1/"

# Row 32: St Josephs 4 - 0 Glacis United
$ws.Range("A32").Value = 31
$ws.Range("B32").Value = "gibraltar"
$ws.Range("C32").Value = "national-league"
$ws.Range("D32").Value = "2023-2024"
$ws.Range("E32").Value = 45255.8125
$ws.Range("F32").Value = "St Josephs"
$ws.Range("G32").Value = 4
$ws.Range("H32").Value = "Glacis United"
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 1.22
$ws.Range("K32").Value = "25/11/2023 12:46"
$ws.Range("L32").Value = 1.08
$ws.Range("M32").Value = "25/11/2023 19:13"
$ws.Range("N32").Value = 6.81
$ws.Range("O32").Value = "25/11/2023 12:46"
$ws.Range("P32").Value = 10.39
$ws.Range("Q32").Value = "25/11/2023 19:22"
$ws.Range("R32").Value = 6.34
$ws.Range("S32").Value = "25/11/2023 12:46"
$ws.Range("T32").Value = 13.06
$ws.Range("U32").Value = "25/11/2023 19:22"
$ws.Range("V32").Value = "https://www.betexplorer.com/football/gibraltar/national-league/st-josephs-glacis-united/67h73Si7/"

Write-Output "Added rows 30-32 to Sheet1"
